$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export Worksheet")
$excel.UserName = "Bruce Higiro Munyandamutsa"
$ws.Columns("P").Insert()
$c1 = $ws.Range("S6").Comment
$txt1 = $c1.Text()
$c1.Delete()
# Leave T6's original comment alone. Add new comment elsewhere (V6):
$new = $ws.Range("V6").AddComment($txt1)
Write-Host "author (live):" $new.Author
Write-Host "T6 original still exists?" ($ws.Range("T6").Comment -ne $null)
